# edit.ps1 - apply JFreeApache Manual edits via Word COM-interop
#
# 1. Title placeholder "XXXXX" -> "JFreeApache" (wrapped in spellcheck proofErr
#    markers, since Word's proofer flags it as an unrecognised word).
# 2. Subtitle placeholder "XXXXXXXXXXXX" -> the project blurb "A modified
#    version of PSS using JFreeCharts and Apache Commons Math" split across
#    several runs (mirrors how Word would have built it up keystroke by
#    keystroke / via autocorrect) with "JFreeCharts" wrapped in proofErr.
# 3. The "Software Description" body paragraph gets centred and its
#    placeholder text replaced with the same project blurb.

$d = $word.ActiveDocument

function Insert-RunsXml($insertionRange, $innerXml) {
    $xml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
$innerXml
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $insertionRange.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 1. Title: "XXXXX" -> "JFreeApache" (proofErr-wrapped), keep the
#    trailing " User Manual" run untouched.
# ---------------------------------------------------------------------------
$pTitle = $d.Paragraphs.Item(6)
$titleText = $pTitle.Range.Text
$xxxxxIdx = $titleText.IndexOf("XXXXX")
if ($xxxxxIdx -ge 0) {
    $titleStart = $pTitle.Range.Start
    $target = $d.Range($titleStart + $xxxxxIdx, $titleStart + $xxxxxIdx + 5)
    $target.Text = ""
    $insertAt = $d.Range($titleStart + $xxxxxIdx, $titleStart + $xxxxxIdx)

    $inner = @"
<w:proofErr w:type="spellStart"/>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
<w:b/>
<w:bCs/>
<w:sz w:val="48"/>
<w:szCs w:val="48"/>
</w:rPr>
<w:t>JFreeApache</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
"@
    Insert-RunsXml $insertAt $inner
}

# ---------------------------------------------------------------------------
# 2. Subtitle: "XXXXXXXXXXXX" -> project blurb, split into several runs,
#    with "JFreeCharts" wrapped in proofErr spellStart/spellEnd.
# ---------------------------------------------------------------------------
$pSubtitle = $d.Paragraphs.Item(7)
$subtitleText = $pSubtitle.Range.Text.TrimEnd("`r", "`n")
$subStart = $pSubtitle.Range.Start
$subLen = $subtitleText.Length
$subTarget = $d.Range($subStart, $subStart + $subLen)
if ($subTarget.Text -eq "XXXXXXXXXXXX") {
    $subTarget.Text = ""
    $subInsertAt = $d.Range($subStart, $subStart)

    $inner = @"
<w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>A modified version of P</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>SS</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> u</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>s</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">ing </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>JFreeChart</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>s</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> and Apache Common</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>s</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Math</w:t></w:r>
"@
    Insert-RunsXml $subInsertAt $inner
}

# ---------------------------------------------------------------------------
# 3. "Software Description" body paragraph: centre it and replace the
#    placeholder sentence with the project blurb (shorter run split, with
#    "JFreeCharts" proofErr-wrapped again).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r", "`n") -eq "Write what the software does here.") {
        $p.Alignment = 1 ## wdAlignParagraphCenter

        $bodyStart = $p.Range.Start
        $bodyLen = $p.Range.Text.TrimEnd("`r", "`n").Length
        $bodyTarget = $d.Range($bodyStart, $bodyStart + $bodyLen)
        $bodyTarget.Text = ""
        $bodyInsertAt = $d.Range($bodyStart, $bodyStart)

        $inner = @"
<w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">A modified version of PSS using </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>JFreeCharts</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> and Apache Commons Math</w:t></w:r>
"@
        Insert-RunsXml $bodyInsertAt $inner
        break
    }
}

Write-Host "Edits applied."
